$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.9999892339301442"
$ws.Range("E2").Value = [double]"0.9999892339301442"

# Row 3
$ws.Range("D3").Value = [double]"0.9999999980687526"
$ws.Range("E3").Value = [double]"0.9999999980687526"

# Row 4
$ws.Range("D4").Value = [double]"2.709214852488288E-05"
$ws.Range("E4").Value = [double]"2.709214852488288E-05"

# Row 5
$ws.Range("D5").Value = [double]"0.04608201996621249"
$ws.Range("E5").Value = [double]"0.04608201996621249"

# Row 6
$ws.Range("D6").Value = [double]"0.01111472498919494"
$ws.Range("E6").Value = [double]"0.01111472498919494"

# Row 7
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = [double]"0.9945513796295906"
$ws.Range("E7").Value = [double]"0.005448620370409385"

# Row 9
$ws.Range("D9").Value = [double]"0.9999999957121937"
$ws.Range("E9").Value = [double]"4.287806332214927E-09"

# Row 10
$ws.Range("D10").Value = [double]"2.37341036269836E-06"
$ws.Range("E10").Value = [double]"0.9999976265896373"

# Row 11
$ws.Range("D11").Value = [double]"1.784662418168611E-05"
$ws.Range("E11").Value = [double]"0.9999821533758183"
$ws.Range("F11").Value = [double]"5.54529333114624"
$ws.Range("G11").Value = [double]"0.6"
